$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column M, row 3: blank cell with the same thin-strip-under-header
# formatting as L3 (font + bottom medium border only). ---
$ws.Range("M3").Font.Name = "Times New Roman"
$ws.Range("M3").Font.Size = 9
$ws.Range("M3").Font.Bold = $false
$ws.Range("M3").VerticalAlignment = -4108
$ws.Range("M3").Borders.Item(9).Weight = -4138
$ws.Range("M3").Borders.Item(9).ColorIndex = 1

# --- New column M, row 4: year header "2022", matching L4's bold/right style. ---
$ws.Range("M4").Font.Name = "Times New Roman"
$ws.Range("M4").Font.Size = 9
$ws.Range("M4").Font.Bold = $true
$ws.Range("M4").HorizontalAlignment = -4152
$ws.Range("M4").VerticalAlignment = -4108
$ws.Range("M4").Borders.Item(9).Weight = -4138
$ws.Range("M4").Borders.Item(9).ColorIndex = 1
$ws.Range("M4").Value = 2022

# --- New column M, row 5: data value "373", matching L5's top+bottom border style. ---
$ws.Range("M5").Font.Name = "Times New Roman"
$ws.Range("M5").Font.Size = 9
$ws.Range("M5").Font.Bold = $false
$ws.Range("M5").VerticalAlignment = -4108
$ws.Range("M5").Borders.Item(8).Weight = -4138
$ws.Range("M5").Borders.Item(8).ColorIndex = 1
$ws.Range("M5").Borders.Item(9).Weight = -4138
$ws.Range("M5").Borders.Item(9).ColorIndex = 1
$ws.Range("M5").Value = 373

# Selection ends up on O4 in the saved file.
$ws.Range("O4").Select()
